# Applies:
#  1. Inserts three new runs ("<comment>", "c_166r_01", "</comment>")
#     immediately before the run that starts "Si le zelotype reprend..."
#  2. Sets the section's footer distance to 720 twips (0.5"), which adds
#     w:footer="720" to the sectPr's pgMar.

$d = $word.ActiveDocument

# --- 1. Insert the <comment>c_166r_01</comment> marker runs -------------

$rng = $d.Content
$rng.Find.Execute("Si le zelotype reprend", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertPoint = $d.Range($rng.Start, $rng.Start)

$xmlSnippet = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/><w:color w:val="0000ff"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">&lt;comment&gt;</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">c_166r_01</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/><w:color w:val="0000ff"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">&lt;/comment&gt;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
# (attribute order for w:rFonts deliberately mirrors the target OOXML:
#  ascii, cs, eastAsia, hAnsi)

$insertPoint.InsertXML($xmlSnippet)

# --- 2. Add w:footer="720" to the section's pgMar ------------------------

$sec = $d.Sections.First
$sec.PageSetup.FooterDistance = 36

Write-Output "edit applied"
